# Update the "Number of APKs" column (column B) on 'Sheet 5' so the chart
# reflects unique APK counts rather than raw counts. Only rows 5-34
# (B5:B34) actually change; rows 3-4 already read 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 5")

$newValues = @(
    1, # B5
    1, # B6
    1, # B7
    2, # B8
    2, # B9
    2, # B10
    2, # B11
    3, # B12
    3, # B13
    3, # B14
    4, # B15
    4, # B16
    6, # B17
    9, # B18
    9, # B19
    10, # B20
    11, # B21
    12, # B22
    12, # B23
    13, # B24
    13, # B25
    14, # B26
    15, # B27
    21, # B28
    27, # B29
    28, # B30
    33, # B31
    33, # B32
    58, # B33
    70  # B34
)

$row = 5
foreach ($v in $newValues) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}
